# 8.5.2.xlsx - add a 2021 data column (column R) to the unemployment-rate
# table on sheet1, mirroring the existing 2007-2020 columns (D:Q).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122
$xlRight = -4152

# ---------------------------------------------------------------------
# 1. Build the two brand-new cell styles that the target workbook needs
#    (in the same order they must appear in cellXfs so the indices line
#    up: new style "17" first, then new style "18").
# ---------------------------------------------------------------------

# Style "17": a fresh Times New Roman / size 9 font (no bold/italic, no
# border, General number format) - used for row 6 and row 7 of the new
# column. Start from the existing plain font-9 style (style 10, e.g.
# cell A8) and nudge the font so the engine forks a brand-new font
# entry instead of reusing font 9.
$ws.Range("A8").Copy()
$ws.Range("R6").PasteSpecial($xlPasteFormats)
$ws.Range("R6").Font.ThemeColor = 1
$ws.Range("R6").Copy()
$ws.Range("R7").PasteSpecial($xlPasteFormats)

# Style "18": same font/border as the existing bottom-border style used
# on row 43 (style 15, e.g. cell A43), plus right alignment - used for
# the new R43 cell.
$ws.Range("A43").Copy()
$ws.Range("R43").PasteSpecial($xlPasteFormats)
$ws.Range("R43").HorizontalAlignment = $xlRight

# ---------------------------------------------------------------------
# 2. Header row: 2021 year label, same style as the other year headers.
# ---------------------------------------------------------------------
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial($xlPasteFormats)
$ws.Range("R4").Value = 2021

# ---------------------------------------------------------------------
# 3. Data rows. Group header rows (both sexes) reuse style 11 (A5-style,
#    bold-free "section" look); the two sub-rows under each group reuse
#    style 10 (A8/A9-style, the plain data-row look) - this mirrors the
#    existing pattern already used in columns A-C.
# ---------------------------------------------------------------------

function Set-R11($row, $value) {
    $ws.Range("A5").Copy()
    $ws.Range("R$row").PasteSpecial($xlPasteFormats)
    $ws.Range("R$row").Value = $value
}

function Set-R10($row, $value) {
    $ws.Range("A8").Copy()
    $ws.Range("R$row").PasteSpecial($xlPasteFormats)
    if ($null -ne $value) {
        $ws.Range("R$row").Value = $value
    }
}

Set-R11 5  5.3
# R6 / R7 already carry the new style "17" from step 1 above.
$ws.Range("R6").Value = 6.3
$ws.Range("R7").Value = 4.7

Set-R10 8  $null

Set-R11 9  6.6
Set-R10 10 7.5
Set-R10 11 6.2

Set-R11 12 11.8
Set-R10 13 15.5
Set-R10 14 9.7

Set-R11 15 6.3
Set-R10 16 7.5
Set-R10 17 5.6

Set-R11 18 6.3
Set-R10 19 10.8
Set-R10 20 4.3

Set-R11 21 1.9
Set-R10 22 3.1
Set-R10 23 1.1

Set-R11 24 2.6
Set-R10 25 3.8
Set-R10 26 1.7

Set-R11 27 5.3
Set-R10 28 6.2
Set-R10 29 4.8

Set-R11 30 4.1
Set-R10 31 3.3
Set-R10 32 4.9

Set-R11 33 2.8
Set-R10 34 3.4
Set-R10 35 2.6

Set-R10 36 $null

Set-R10 37 15.7
Set-R10 38 7.9
Set-R10 39 4.5
Set-R10 40 4.4
Set-R10 41 2.9
Set-R10 42 1.4

# ---------------------------------------------------------------------
# 4. Footer row: "..." placeholder, same shared string already used
#    elsewhere in the sheet, with the new right-aligned bottom-border
#    style "18" created in step 1.
# ---------------------------------------------------------------------
$ws.Range("R43").Value = [char]0x2026

# ---------------------------------------------------------------------
# 5. Selection / active cell, matching the post-edit state recorded in
#    the workbook.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("S1").Select()
